$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.849.35'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").Value = '2.529.66'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '317.87'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '96.73'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = '35.85'
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").Value = '7.52'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("D14").Value = '2.918.51'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '2.561.41'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").Value = '15.06'
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").Value = '42.867.41'
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("E19").Value = '  +2.93%  '
$ws.Range("D20").Value = '12.56'
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").Value = '0.0₃0965'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("D22").Value = '69.60'
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("D23").Value = '252.84'
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '2.05'
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").Value = '26.38'
$ws.Range("E26").Value = '  -4.72%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").Value = '41.35'
$ws.Range("E29").Value = '  +4.60%  '
$ws.Range("D30").Value = '10.42'
$ws.Range("E30").Value = '  +3.32%  '
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").Value = '156.90'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = '2.14'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '19.38'
$ws.Range("E34").Value = '  -0.60%  '
$ws.Range("D35").Value = '3.36'
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("E39").Value = '  +9.16%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '21.90'
$ws.Range("E41").Value = '  -11.79%  '
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = '3.80'
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("E44").Value = '  +0.27%  '
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("D46").Value = '1.998.25'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("D48").Value = '84.40'
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '105.85'
$ws.Range("E49").Value = '  +3.42%  '
$ws.Range("D50").Value = '75.03'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '2.773.05'
$ws.Range("E51").Value = '  -0.47%  '
